$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "Текст до внесения изменений:" + break + "1" + break + "2" + break + "3" + break + "4" + break + "5"
# It needs to become:
#   "Текст " / "п" / "о" / "сле" / " внесения изменений:"  (5 separate runs spelling
#   "Текст после внесения изменений:") followed by a single, now-empty, line break run.

$p = $d.Paragraphs(1)

# Remember where the paragraph starts so we can compute the boundary between the
# newly-inserted text and the stale content that must be removed.
$paraStart = $p.Range.Start

# Collapse a range to the very start of the paragraph; repeated InsertBefore calls on a
# range collapsed this way each land immediately at that fixed point, so inserting the
# pieces in reverse order produces them, left to right, as separate runs without merging
# into the (still present) original runs that follow.
$ip = $p.Range
$ip.Collapse(1)

$ip.InsertBreak(6)
$ip.InsertBefore(" внесения изменений:")
$ip.InsertBefore("сле")
$ip.InsertBefore("о")
$ip.InsertBefore("п")
$ip.InsertBefore("Текст ")

$newLen = ("Текст после внесения изменений:").Length + 1
$tailStart = $paraStart + $newLen

$paraEnd = $p.Range.End
$tail = $d.Range($tailStart, $paraEnd - 1)
$tail.Delete()
